$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.125615
$ws.Range("H2").Value = 0.25123
$ws.Range("I2").Value = 0.02647478672532295
$ws.Range("J2").Value = 0.01780700335556722
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.532593
$ws.Range("N2").Value = 3.065186
$ws.Range("O2").Value = 0.0795983245703594
$ws.Range("P2").Value = 0.05838920196386116
$ws.Range("Q2").Value = 0.192516669695
$ws.Range("R2").Value = 0.7700666787800001
$ws.Range("S2").Value = 0.002107348666693299
$ws.Range("T2").Value = 0.001039736715299368

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.125615
$ws.Range("H3").Value = 0.25123
$ws.Range("I3").Value = 0.02647478672532295
$ws.Range("J3").Value = 0.01780700335556722
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.18022033333333
$ws.Range("N3").Value = 36.540661
$ws.Range("O3").Value = 0.6326044366842063
$ws.Range("P3").Value = 0.6960687002426557
$ws.Range("Q3").Value = 1.530018377171667
$ws.Range("R3").Value = 9.18011026303
$ws.Range("S3").Value = 0.01674806754270743
$ws.Range("T3").Value = 0.01239489768092628

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.125615
$ws.Range("H4").Value = 0.25123
$ws.Range("I4").Value = 0.02647478672532295
$ws.Range("J4").Value = 0.01780700335556722
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8528209999999999
$ws.Range("N4").Value = 2.558463
$ws.Range("O4").Value = 0.04429298760885536
$ws.Range("P4").Value = 0.04873655720209673
$ws.Range("Q4").Value = 0.107127109915
$ws.Range("R4").Value = 0.64276265949
$ws.Range("S4").Value = 0.001172647400371818
$ws.Range("T4").Value = 0.0008678520376365303

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.125615
$ws.Range("H5").Value = 0.25123
$ws.Range("I5").Value = 0.02647478672532295
$ws.Range("J5").Value = 0.01780700335556722
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.7338975
$ws.Range("N5").Value = 7.467795
$ws.Range("O5").Value = 0.1939275366111247
$ws.Range("P5").Value = 0.142255181408147
$ws.Range("Q5").Value = 0.4690335344625
$ws.Range("R5").Value = 1.87613413785
$ws.Range("S5").Value = 0.005134190171946787
$ws.Range("T5").Value = 0.002533138492681698

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd1"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.125615
$ws.Range("H6").Value = 0.25123
$ws.Range("I6").Value = 0.02647478672532295
$ws.Range("J6").Value = 0.01780700335556722
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2147316666666667
$ws.Range("N6").Value = 0.644195
$ws.Range("O6").Value = 0.01115252444639089
$ws.Range("P6").Value = 0.01227137014168456
$ws.Range("Q6").Value = 0.02697351830833333
$ws.Range("R6").Value = 0.16184110985
$ws.Range("S6").Value = 0.0002952607061671493
$ws.Range("T6").Value = 0.0002185163292903844

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.125615
$ws.Range("H7").Value = 0.25123
$ws.Range("I7").Value = 0.02647478672532295
$ws.Range("J7").Value = 0.01780700335556722
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7398226666666666
$ws.Range("N7").Value = 2.219468
$ws.Range("O7").Value = 0.03842419007906348
$ws.Range("P7").Value = 0.04227898904155473
$ws.Range("Q7").Value = 0.09293282427333333
$ws.Range("R7").Value = 0.5575969456400001
$ws.Range("S7").Value = 0.001017272237436476
$ws.Range("T7").Value = 0.0007528620997329548

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.619088000000001
$ws.Range("H8").Value = 13.857264
$ws.Range("I8").Value = 0.9735252132746771
$ws.Range("J8").Value = 0.9821929966444328
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.532593
$ws.Range("N8").Value = 3.065186
$ws.Range("O8").Value = 0.0795983245703594
$ws.Range("P8").Value = 0.05838920196386116
$ws.Range("Q8").Value = 7.079181935184002
$ws.Range("R8").Value = 42.47509161110401
$ws.Range("S8").Value = 0.0774909759036661
$ws.Range("T8").Value = 0.05734946524856179

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.619088000000001
$ws.Range("H9").Value = 13.857264
$ws.Range("I9").Value = 0.9735252132746771
$ws.Range("J9").Value = 0.9821929966444328
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.18022033333333
$ws.Range("N9").Value = 36.540661
$ws.Range("O9").Value = 0.6326044366842063
$ws.Range("P9").Value = 0.6960687002426557
$ws.Range("Q9").Value = 56.261509579056
$ws.Range("R9").Value = 506.353586211504
$ws.Range("S9").Value = 0.6158563691414989
$ws.Range("T9").Value = 0.6836738025617294

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt5a"
$ws.Range("C10").Value = "Fzd1"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.619088000000001
$ws.Range("H10").Value = 13.857264
$ws.Range("I10").Value = 0.9735252132746771
$ws.Range("J10").Value = 0.9821929966444328
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8528209999999999
$ws.Range("N10").Value = 2.558463
$ws.Range("O10").Value = 0.04429298760885536
$ws.Range("P10").Value = 0.04873655720209673
$ws.Range("Q10").Value = 3.939255247248
$ws.Range("R10").Value = 35.453297225232
$ws.Range("S10").Value = 0.04312034020848354
$ws.Range("T10").Value = 0.04786870516446021

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt5a"
$ws.Range("C11").Value = "Fzd1"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.619088000000001
$ws.Range("H11").Value = 13.857264
$ws.Range("I11").Value = 0.9735252132746771
$ws.Range("J11").Value = 0.9821929966444328
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.7338975
$ws.Range("N11").Value = 7.467795
$ws.Range("O11").Value = 0.1939275366111247
$ws.Range("P11").Value = 0.142255181408147
$ws.Range("Q11").Value = 17.24720113548
$ws.Range("R11").Value = 103.48320681288
$ws.Range("S11").Value = 0.188793346439178
$ws.Range("T11").Value = 0.1397220429154653

# Row 12
$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Wnt5a"
$ws.Range("C12").Value = "Fzd1"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.619088000000001
$ws.Range("H12").Value = 13.857264
$ws.Range("I12").Value = 0.9735252132746771
$ws.Range("J12").Value = 0.9821929966444328
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2147316666666667
$ws.Range("N12").Value = 0.644195
$ws.Range("O12").Value = 0.01115252444639089
$ws.Range("P12").Value = 0.01227137014168456
$ws.Range("Q12").Value = 0.99186446472
$ws.Range("R12").Value = 8.92678018248
$ws.Range("S12").Value = 0.01085726374022374
$ws.Range("T12").Value = 0.01205285381239418

# Row 13
$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Wnt5a"
$ws.Range("C13").Value = "Fzd1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.619088000000001
$ws.Range("H13").Value = 13.857264
$ws.Range("I13").Value = 0.9735252132746771
$ws.Range("J13").Value = 0.9821929966444328
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7398226666666666
$ws.Range("N13").Value = 2.219468
$ws.Range("O13").Value = 0.03842419007906348
$ws.Range("P13").Value = 0.04227898904155473
$ws.Range("Q13").Value = 3.417306001728
$ws.Range("R13").Value = 30.755754015552
$ws.Range("S13").Value = 0.03740691784162701
$ws.Range("T13").Value = 0.04152612694182178
